$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 34.42857
$ws.Range("I11").Value = 34.42857
$ws.Range("K11").Value = 34.42857
$ws.Range("M11").Value = 105.57143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2975
$ws.Range("J51").Value = 2975
$ws.Range("L51").Value = 2975
$ws.Range("N51").Value = -3943

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1499.5
$ws.Range("I100").Value = 999
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 999
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -458
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7873
$ws.Range("I116").Value = 8716.25
$ws.Range("K116").Value = 8716.25
$ws.Range("M116").Value = -5274.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14793.066
$ws.Range("I132").Value = 15678.643
$ws.Range("K132").Value = 47035.929
$ws.Range("M132").Value = -44505.929

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 934
$ws.Range("I135").Value = 907.06665
$ws.Range("J135").Value = 991.7143
$ws.Range("K135").Value = 8163.59985
$ws.Range("L135").Value = 8925.4287
$ws.Range("M135").Value = -5628.59985
$ws.Range("N135").Value = -13995.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3971.1904
$ws.Range("I137").Value = 999.2308
$ws.Range("J137").Value = 8800.625
$ws.Range("K137").Value = 2997.6924
$ws.Range("L137").Value = 26401.875
$ws.Range("M137").Value = -447.6923999999999
$ws.Range("N137").Value = -31501.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4498.8335
$ws.Range("I138").Value = 1448.25
$ws.Range("J138").Value = 10600
$ws.Range("K138").Value = 4344.75
$ws.Range("L138").Value = 31800
$ws.Range("M138").Value = 795.25
$ws.Range("N138").Value = -42080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1081.5
$ws.Range("I2").Value = 1081.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1081.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -968.5
$ws.Range("N2").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2039.8
$ws.Range("I26").Value = 2039.8
$ws.Range("K26").Value = 2039.8
$ws.Range("M26").Value = -1709.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 24000
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2967.5715
$ws.Range("I45").Value = 2170.5
$ws.Range("J45").Value = 3286.4
$ws.Range("K45").Value = 2170.5
$ws.Range("L45").Value = 3286.4
$ws.Range("M45").Value = -1793.5
$ws.Range("N45").Value = -4040.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 982.3333
$ws.Range("I88").Value = 404.16666
$ws.Range("J88").Value = 2138.6667
$ws.Range("K88").Value = 404.16666
$ws.Range("L88").Value = 2138.6667
$ws.Range("M88").Value = 1.833340000000021
$ws.Range("N88").Value = -2950.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 982.3333
$ws.Range("I91").Value = 404.16666
$ws.Range("J91").Value = 2138.6667
$ws.Range("K91").Value = 404.16666
$ws.Range("L91").Value = 2138.6667
$ws.Range("M91").Value = 999.83334
$ws.Range("N91").Value = -4946.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 965.5
$ws.Range("I97").Value = 928.8
$ws.Range("J97").Value = 1149
$ws.Range("K97").Value = 928.8
$ws.Range("L97").Value = 1149
$ws.Range("M97").Value = -432.8
$ws.Range("N97").Value = -2141

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 17587.25
$ws.Range("J112").Value = 17587.25
$ws.Range("L112").Value = 17587.25
$ws.Range("N112").Value = -20541.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1081.5
$ws.Range("I116").Value = 1081.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1081.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1212.5
$ws.Range("N116").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2416.5
$ws.Range("I122").Value = 2416.5
$ws.Range("K122").Value = 7249.5
$ws.Range("M122").Value = -4799.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1961.6
$ws.Range("I132").Value = 1670.1052
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 5010.3156
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -2480.3156
$ws.Range("N132").Value = -27560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1081.5
$ws.Range("I3").Value = 1081.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1081.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -967.5
$ws.Range("N3").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 13341844
$ws.Range("I7").Value = 16000255
$ws.Range("J7").Value = 8025020
$ws.Range("K7").Value = 16000255
$ws.Range("L7").Value = 8025020
$ws.Range("M7").Value = -16000142
$ws.Range("N7").Value = -8025246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 673.4
$ws.Range("I94").Value = 667
$ws.Range("K94").Value = 667
$ws.Range("M94").Value = -216

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2637.5
$ws.Range("I105").Value = 2275
$ws.Range("K105").Value = 2275
$ws.Range("M105").Value = -528

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1441.7
$ws.Range("I35").Value = 1052.25
$ws.Range("K35").Value = 1052.25
$ws.Range("M35").Value = -758.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3281
$ws.Range("I58").Value = 1915
$ws.Range("K58").Value = 1915
$ws.Range("M58").Value = -1712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 825
$ws.Range("I105").Value = 766.6667
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 766.6667
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 980.3333
$ws.Range("N105").Value = -4494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 798.25
$ws.Range("I107").Value = 497
$ws.Range("J107").Value = 1013.4286
$ws.Range("K107").Value = 497
$ws.Range("L107").Value = 1013.4286
$ws.Range("M107").Value = 1423
$ws.Range("N107").Value = -4853.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2984.8147
$ws.Range("I132").Value = 2967.6
$ws.Range("K132").Value = 8902.799999999999
$ws.Range("M132").Value = -6372.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1644.0667
$ws.Range("I134").Value = 973.88464
$ws.Range("J134").Value = 6000.25
$ws.Range("K134").Value = 2921.65392
$ws.Range("L134").Value = 18000.75
$ws.Range("M134").Value = -386.6539199999997
$ws.Range("N134").Value = -23070.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3281
$ws.Range("I136").Value = 1915
$ws.Range("K136").Value = 5745
$ws.Range("M136").Value = -3195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 1425
$ws.Range("I124").Value = 1425
$ws.Range("K124").Value = 4275
$ws.Range("M124").Value = 635

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3261.7273
$ws.Range("I140").Value = 2819.889
$ws.Range("K140").Value = 8459.667000000001
$ws.Range("M140").Value = -3279.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2225.5334
$ws.Range("I102").Value = 2241.6428
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2241.6428
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -619.6428000000001
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 44000
$ws.Range("J110").Value = 44000
$ws.Range("L110").Value = 44000
$ws.Range("N110").Value = -52180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 37157.3
$ws.Range("I132").Value = 45351.668
$ws.Range("K132").Value = 136055.004
$ws.Range("M132").Value = -133525.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1693
$ws.Range("I10").Value = 1252
$ws.Range("J10").Value = 2575
$ws.Range("K10").Value = 1252
$ws.Range("L10").Value = 2575
$ws.Range("M10").Value = -1112
$ws.Range("N10").Value = -2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2179.6365
$ws.Range("I22").Value = 850
$ws.Range("J22").Value = 3287.6667
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 3287.6667
$ws.Range("M22").Value = -555
$ws.Range("N22").Value = -3877.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2179.6365
$ws.Range("I27").Value = 850
$ws.Range("J27").Value = 3287.6667
$ws.Range("K27").Value = 850
$ws.Range("L27").Value = 3287.6667
$ws.Range("M27").Value = -743
$ws.Range("N27").Value = -3501.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3724.875
$ws.Range("I61").Value = 2050
$ws.Range("K61").Value = 2050
$ws.Range("M61").Value = -1848

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3724.875
$ws.Range("I113").Value = 2050
$ws.Range("K113").Value = 2050
$ws.Range("M113").Value = 120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3370.3572
$ws.Range("I122").Value = 3456
$ws.Range("J122").Value = 3284.7144
$ws.Range("K122").Value = 10368
$ws.Range("L122").Value = 9854.143199999999
$ws.Range("M122").Value = -7918
$ws.Range("N122").Value = -14754.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2670.842
$ws.Range("I122").Value = 2174.6924
$ws.Range("J122").Value = 3745.8333
$ws.Range("K122").Value = 6524.0772
$ws.Range("L122").Value = 11237.4999
$ws.Range("M122").Value = -4074.0772
$ws.Range("N122").Value = -16137.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2502.1365
$ws.Range("I126").Value = 2129.2942
$ws.Range("K126").Value = 6387.882599999999
$ws.Range("M126").Value = -3917.882599999999
